# Apply price/volume/coin-ranking updates scraped on 2023-08-10 (GitHub Actions job).
# Mirrors the authoritative OOXML diff cell-for-cell:
#   - D/E columns on rows 2-51 get refreshed Price / Volume(1h) figures.
#   - Rows 45 and 46 swap rank (Aave <-> Quant) together with their data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking Price strings (single "." or none) need a leading apostrophe so
# Excel keeps them as literal text instead of inferring a Number/Date value -
# exactly like the source file, which stores every cell in this table as a string.
# Price strings that already contain two "." (thousands-separator style, e.g. 29.404.91)
# are never auto-converted, so they are assigned as-is.

# Row 2
$ws.Range("D2").Value = "29.404.91"
$ws.Range("E2").Value = "  -0.54%  "
# Row 3
$ws.Range("D3").Value = "1.849.70"
$ws.Range("E3").Value = "  -0.09%  "
# Row 4
$ws.Range("D4").Value = "'0.9992"
$ws.Range("E4").Value = "  -0.09%  "
# Row 5
$ws.Range("D5").Value = "'241.19"
$ws.Range("E5").Value = "  -0.93%  "
# Row 6
$ws.Range("D6").Value = "'0.6265"
$ws.Range("E6").Value = "  -3.87%  "
# Row 7
$ws.Range("E7").Value = "  -0.06%  "
# Row 8
$ws.Range("D8").Value = "'0.07635"
$ws.Range("E8").Value = "  +1.90%  "
# Row 9
$ws.Range("D9").Value = "'0.2971"
$ws.Range("E9").Value = "  -0.40%  "
# Row 10
$ws.Range("D10").Value = "'24.45"
$ws.Range("E10").Value = "  -0.38%  "
# Row 11
$ws.Range("D11").Value = "2.103.31"
$ws.Range("E11").Value = "  +13.62%  "
# Row 12
$ws.Range("D12").Value = "'0.07726"
$ws.Range("E12").Value = "  +1.09%  "
# Row 13
$ws.Range("D13").Value = "'4.997"
$ws.Range("E13").Value = "  -0.87%  "
# Row 14
$ws.Range("D14").Value = "'0.6881"
$ws.Range("E14").Value = "  +0.10%  "
# Row 15
$ws.Range("D15").Value = "'83.00"
$ws.Range("E15").Value = "  -0.71%  "
# Row 16
$ws.Range("D16").Value = "'0.000009953"
$ws.Range("E16").Value = "  +4.36%  "
# Row 17
$ws.Range("D17").Value = "2.268.92"
$ws.Range("E17").Value = "  +7.91%  "
# Row 18
$ws.Range("D18").Value = "'6.165"
$ws.Range("E18").Value = "  +0.67%  "
# Row 19
$ws.Range("D19").Value = "29.649.47"
$ws.Range("E19").Value = "  +0.17%  "
# Row 20
$ws.Range("D20").Value = "'231.24"
$ws.Range("E20").Value = "  -2.59%  "
# Row 21
$ws.Range("D21").Value = "'12.54"
$ws.Range("E21").Value = "  -0.63%  "
# Row 22
$ws.Range("D22").Value = "'1.000"
$ws.Range("E22").Value = "  -0.05%  "
# Row 23
$ws.Range("D23").Value = "'7.643"
$ws.Range("E23").Value = "  -1.18%  "
# Row 24
$ws.Range("E24").Value = "  -0.16%  "
# Row 25
$ws.Range("D25").Value = "'154.78"
$ws.Range("E25").Value = "  -1.81%  "
# Row 26
$ws.Range("D26").Value = "'0.1389"
$ws.Range("E26").Value = "  -2.39%  "
# Row 27
$ws.Range("D27").Value = "'8.480"
$ws.Range("E27").Value = "  -0.58%  "
# Row 28
$ws.Range("D28").Value = "'17.67"
$ws.Range("E28").Value = "  -0.98%  "
# Row 29
$ws.Range("E29").Value = "  -0.94%  "
# Row 30
$ws.Range("D30").Value = "'0.05795"
$ws.Range("E30").Value = "  -4.77%  "
# Row 31
$ws.Range("D31").Value = "'1.257"
$ws.Range("E31").Value = "  -0.27%  "
# Row 32
$ws.Range("D32").Value = "'4.129"
$ws.Range("E32").Value = "  -0.31%  "
# Row 33
$ws.Range("D33").Value = "'4.018"
$ws.Range("E33").Value = "  -1.48%  "
# Row 34
$ws.Range("D34").Value = "'1.874"
$ws.Range("E34").Value = "  +0.23%  "
# Row 35
$ws.Range("D35").Value = "'1.161"
$ws.Range("E35").Value = "  -2.69%  "
# Row 36
$ws.Range("D36").Value = "'0.7178"
$ws.Range("E36").Value = "  -1.39%  "
# Row 37
$ws.Range("D37").Value = "'2.589"
$ws.Range("E37").Value = "  -0.42%  "
# Row 38
$ws.Range("D38").Value = "1.255.51"
$ws.Range("E38").Value = "  +4.49%  "
# Row 39
$ws.Range("D39").Value = "'2.795"
$ws.Range("E39").Value = "  -0.32%  "
# Row 40
$ws.Range("D40").Value = "'0.01804"
$ws.Range("E40").Value = "  +1.08%  "
# Row 41
$ws.Range("D41").Value = "'0.9083"
$ws.Range("E41").Value = "  -0.14%  "
# Row 42
$ws.Range("D42").Value = "2.179.31"
$ws.Range("E42").Value = "  +8.41%  "
# Row 43
$ws.Range("D43").Value = "'6.081"
$ws.Range("E43").Value = "  -2.32%  "
# Row 44
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  -0.08%  "
# Row 45
$ws.Range("B45").Value = "Quant"
$ws.Range("C45").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D45").Value = "'101.91"
$ws.Range("E45").Value = "  +0.38%  "
# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").Value = "'67.53"
$ws.Range("E46").Value = "  +1.39%  "
# Row 47
$ws.Range("D47").Value = "'7.314"
$ws.Range("E47").Value = "  -0.87%  "
# Row 48
$ws.Range("E48").Value = "  -3.34%  "
# Row 49
$ws.Range("D49").Value = "'9.171"
$ws.Range("E49").Value = "  +0.52%  "
# Row 50
$ws.Range("D50").Value = "'0.4030"
$ws.Range("E50").Value = "  -0.68%  "
# Row 51
$ws.Range("E51").Value = "  +1.80%  "
